$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 116 currently holds the "old" record for Berenjena at Macroferia Regional
# de Talca. The edit moves that old record down to a newly-appended row 117,
# and updates row 116 in place with the new weekly figures.

# 1) Copy the current (old) row 116 values down into the new row 117.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(117, $col).Value2 = $ws.Cells.Item(116, $col).Value2
}
# Match the date cell's number format (style) used on column D.
$ws.Cells.Item(117, 4).NumberFormat = $ws.Cells.Item(116, 4).NumberFormat

# 2) Update row 116 in place with the new weekly figures.
$ws.Cells.Item(116, 4).Value2 = 44628   # Fecha
$ws.Cells.Item(116, 10).Value2 = 200    # Volumen
$ws.Cells.Item(116, 11).Value2 = 7000   # Precio minimo
$ws.Cells.Item(116, 12).Value2 = 7000   # Precio maximo
$ws.Cells.Item(116, 13).Value2 = 7000   # Precio promedio ponderado
$ws.Cells.Item(116, 16).Value2 = 140    # Precio $/Kg
